$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set up the additional column widths used by the new B/C/D columns of data
# (stored widths are pixel-quantized by the engine; these inputs give the
# closest achievable match to the target widths)
$ws.Columns.Item(2).ColumnWidth = 17.35
$ws.Columns.Item(3).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 16.85

# --- Block 1: "c = 0.1, 16 samples vs 16 samples" (rows 1-6) ---
$ws.Range("A1").Value = "c = 0.1, 16 samples vs 16 samples"

$ws.Range("A2").Value = "Linear"
$ws.Range("B2").Value = "Poly"
$ws.Range("C2").Value = "RBF"
$ws.Range("D2").Value = "sigmoid"

$ws.Range("A3").Value = "SVMAccuracy : 78.125"
$ws.Range("B3").Value = "SVMAccuracy : 62.5"
$ws.Range("C3").Value = "SVMAccuracy : 62.5"
$ws.Range("D3").Value = "SVMAccuracy : 68.75"

$ws.Range("A4").Value = "SVMConfusionMatrix"
$ws.Range("B4").Value = "SVMConfusionMatrix"
$ws.Range("C4").Value = "SVMConfusionMatrix"
$ws.Range("D4").Value = "SVMConfusionMatrix"

$ws.Range("A5").Value = "81  19"
$ws.Range("B5").Value = "50  50"
$ws.Range("C5").Value = "100    0"
$ws.Range("D5").Value = "75  25"

$ws.Range("A6").Value = "25  75"
$ws.Range("B6").Value = "25  75"
$ws.Range("C6").Value = " 75   25"
$ws.Range("D6").Value = "38  63"

# --- Block 2: "c = 0.1, half and half, SVM vs MKL" (rows 8-16) ---
$ws.Range("A8").Value = "c = 0.1, half and half, SVM vs MKL"

$ws.Range("A9").Value = "SVMAccuracy : 75"
$ws.Range("B9").Value = "SVMAccuracy : 50"
$ws.Range("C9").Value = "SVMAccuracy : 62.5"
$ws.Range("D9").Value = "SVMAccuracy : 87.5"

$ws.Range("A10").Value = "MKLAccuracy : 68.75"
$ws.Range("B10").Value = "MKLAccuracy : 68.75"
$ws.Range("C10").Value = "MKLAccuracy : 68.75"
$ws.Range("D10").Value = "MKLAccuracy : 68.75"

$ws.Range("A11").Value = "SVMConfusionMatrix"
$ws.Range("B11").Value = "SVMConfusionMatrix"
$ws.Range("C11").Value = "SVMConfusionMatrix"
$ws.Range("D11").Value = "SVMConfusionMatrix"

$ws.Range("A12").Value = "63  38"
$ws.Range("B12").Value = "50  50"
$ws.Range("C12").Value = "100    0"
$ws.Range("D12").Value = "88  13"

$ws.Range("A13").Value = "13  88"
$ws.Range("B13").Value = "50  50"
$ws.Range("C13").Value = " 75   25"
$ws.Range("D13").Value = "13  88"

$ws.Range("A14").Value = "MKLConfusionMatrix"
$ws.Range("B14").Value = "MKLConfusionMatrix"
$ws.Range("C14").Value = "MKLConfusionMatrix"
$ws.Range("D14").Value = "MKLConfusionMatrix"

$ws.Range("A15").Value = "63  38"
$ws.Range("B15").Value = "63  38"
$ws.Range("C15").Value = "63  38"
$ws.Range("D15").Value = "63  38"

$ws.Range("A16").Value = "25  75"
$ws.Range("B16").Value = "25  75"
$ws.Range("C16").Value = "25  75"
$ws.Range("D16").Value = "25  75"

# Match the saved selection/active cell from the authored workbook
$ws.Range("D24").Select() | Out-Null
